# Update the "想去人数" (want-to-go count) figures in column F across the
# "展览" (Exhibition), "本地生活" (Local Life) and "全部类型" (All Types) sheets,
# reflecting a refreshed scrape (gh-pages output regenerated at 456a3b4).

$wb = $excel.ActiveWorkbook

$wsExpo = $wb.Worksheets.Item("展览")
$wsExpo.Range("F6").Value  = 3821
$wsExpo.Range("F10").Value = 3090
$wsExpo.Range("F13").Value = 2297
$wsExpo.Range("F23").Value = 364
$wsExpo.Range("F33").Value = 4258
$wsExpo.Range("F34").Value = 3945
$wsExpo.Range("F40").Value = 465
$wsExpo.Range("F48").Value = 58

$wsLocal = $wb.Worksheets.Item("本地生活")
$wsLocal.Range("F4").Value = 2273

$wsAll = $wb.Worksheets.Item("全部类型")
$wsAll.Range("F10").Value = 3821
$wsAll.Range("F14").Value = 3090
$wsAll.Range("F16").Value = 2297
$wsAll.Range("F24").Value = 364
$wsAll.Range("F35").Value = 4258
$wsAll.Range("F36").Value = 3945
$wsAll.Range("F40").Value = 465
$wsAll.Range("F48").Value = 58
